# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.658.80"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.730.42"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.82"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.90"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").Value = "3.728.76"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.57"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.479"
$ws.Range("E12").Value = "  -4.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.73"
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000253"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "4.349.89"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "3.726.00"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "69.680.69"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "501.53"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.32"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.14"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.59"
$ws.Range("E24").Value = "  +4.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.16"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.31"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.91"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000136"
$ws.Range("E28").Value = "  +6.91%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.46"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.91"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.32"
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.09"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.349"
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("E39").Value = "  +4.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.03"
$ws.Range("E40").Value = "  +11.25%  "
$ws.Range("B41").Value = "Arweave"
$ws.Range("C41").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "45.81"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.05"
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.64"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "436.89"
$ws.Range("E44").Value = "  +3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.56"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("D46").Value = "2.953.33"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0361"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.36"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.09"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  -2.53%  "

Write-Output "Applied 96 cell updates (35 text-protected)"
